$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.657.49"
$ws.Range("D3").Value = "2.978.99"
$ws.Range("E3").Value = "  -5.05%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.70%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").Value = "2.991.75"
$ws.Range("E9").Value = "  -4.97%  "
$ws.Range("E10").Value = "  -3.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.01%  "
$ws.Range("E12").Value = "  -3.92%  "
$ws.Range("D13").Value = "3.502.01"
$ws.Range("E13").Value = "  -4.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.124"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").Value = "61.717.31"
$ws.Range("E15").Value = "  -4.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.64%  "
$ws.Range("D17").Value = "2.981.12"
$ws.Range("E17").Value = "  -5.15%  "
$ws.Range("E18").Value = "  -5.53%  "
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("E20").Value = "  -3.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.42%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  -3.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.45%  "
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("D27").Value = "3.101.79"
$ws.Range("E27").Value = "  -5.28%  "
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "0.0₃0938"
$ws.Range("E30").Value = "  -7.89%  "
$ws.Range("E31").Value = "  -7.84%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  -4.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.66"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("E36").Value = "  -5.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.04%  "
$ws.Range("E38").Value = "  -4.79%  "
$ws.Range("E39").Value = "  -6.18%  "
$ws.Range("E40").Value = "  -8.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("D42").Value = "2.419.38"
$ws.Range("E42").Value = "  -8.24%  "
$ws.Range("E43").Value = "  -4.40%  "
$ws.Range("E44").Value = "  -6.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.672"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0591"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.44%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("E49").Value = "  -3.85%  "
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("E51").Value = "  -6.82%  "